$d = $word.ActiveDocument

function New-WordXmlFragment([string]$innerWordML) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' + $innerWordML + '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# 1) Insert a brand-new paragraph at the very start of the document containing "{title}"
$startRange = $d.Range(0, 0)
$startRange.InsertParagraphBefore()
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Text = "{title}"

# 2) "Soru {index}. {question}" -> "Question  {index}. {question}" (with proofing marks)
$questionRange = $d.Content
[void]$questionRange.Find.Execute("Soru {index}. {question}")
$questionPara = $questionRange.Paragraphs(1)
$questionParaRange = $questionPara.Range
$questionFull = $d.Range($questionParaRange.Start, $questionParaRange.End - 1)
$questionXml = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Question</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> {</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>index</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}. {</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>question</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}</w:t></w:r>'
$questionFull.InsertXML((New-WordXmlFragment $questionXml))

# 3) "✔ Cevaplar ve Açıklamalar:" -> "✔ Answers and Explanations" (with proofing marks); keep the first (symbol) run untouched
$symbolRange = $d.Content
[void]$symbolRange.Find.Execute("✔")
$afterSymbolStart = $symbolRange.End
$headerPara = $symbolRange.Paragraphs(1)
$headerParaRange = $headerPara.Range
$restOfHeader = $d.Range($afterSymbolStart, $headerParaRange.End - 1)
$answersXml = '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Answers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Explanations</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$restOfHeader.InsertXML((New-WordXmlFragment $answersXml))

# 4) "{index}. {answer} - {explanation}" -> same text, but with split runs + gramStart/gramEnd proofing marks
$answerRange = $d.Content
[void]$answerRange.Find.Execute("{index}. {answer} - {explanation}")
$answerPara = $answerRange.Paragraphs(1)
$answerParaRange = $answerPara.Range
$answerFull = $d.Range($answerParaRange.Start, $answerParaRange.End - 1)
$answerXml = '<w:r><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>index</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}. {</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>answer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>} -</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> {</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>explanation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r>'
$answerFull.InsertXML((New-WordXmlFragment $answerXml))
